$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Layer0")
$ws1.Range("B2").Value = -0.01257756787104805
$ws1.Range("C2").Value = -0.9298945772194686
$ws1.Range("B3").Value = -0.3357097896181181
$ws1.Range("C3").Value = 0.1237095343932099
$ws1.Range("B4").Value = -1.142490738400083
$ws1.Range("C4").Value = -0.2269153530895137

$ws2 = $wb.Worksheets.Item("Layer1")
$ws2.Range("B2").Value = -0.7106686183857307
$ws2.Range("C2").Value = -0.4309124921599725
$ws2.Range("B3").Value = -1.267004070887547
$ws2.Range("C3").Value = 0.6290271227358801
$ws2.Range("B4").Value = 0.1126932270918542
$ws2.Range("C4").Value = 0.7808488167175137
